$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the k column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style B14: bold, size 12, vertical-centered -- then propagate the same
# formatting (without disturbing the formulas) to B15:B17 via copy/paste of
# formats only, so a single new cell style is produced instead of several
# transient ones.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights for the new summary rows
$ws.Range("A14:B17").RowHeight = 15.6

# Page setup (A4, portrait) matching the authored print settings
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection matches the authored sheet view
$ws.Range("A14:B17").Select()

$wb.Save()
